$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J. Everything from the old column J
# onward (ExpectedFilenames, StudyDesignExpectedValue, ReportedVarExpectedValue, ...)
# shifts right by one column (J->K, K->L, L->M, ...).
$ws.Columns("J").Insert()

# The old "ExpectedSourceTemplateFile" header/column is now split into two:
# column I becomes the Excel-report expected template file, and the newly
# inserted column J becomes the Word-report expected template file.
# Set J1 first so the shared-string table allocates the "_Word" text before
# the "_Excel" text (matches the order produced by the original edit).
$ws.Range("J1").Value = "ExpectedSourceTemplateFile_Word"
$ws.Range("I1").Value = "ExpectedSourceTemplateFile_Excel"

# Populate the new column J (Word-report expected template paths) for the
# first few data rows.
$ws.Range("J2").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\Clinical.xlsx"
$ws.Range("J3").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\Economic.xlsx"
$ws.Range("J4").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\QOL.xlsx"
$ws.Range("J5").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\RWE.xlsx"

# Give the new column the same width as column I (the column it was split from).
$ws.Columns("J").ColumnWidth = $ws.Columns("I").ColumnWidth

# Update the view: scroll a bit to the right and move the active selection.
$ws.Range("K5").Select()
